$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.460.97"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.646.35"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.57%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "299.10"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3785"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.95%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3532"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "49.99"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08081"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.213"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.15%  "

$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.08"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.324"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.71%  "

$ws.Range("E16").Value = "  -3.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.644.01"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.97"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06964"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.742"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.40"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.41"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.478.08"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.501"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.879"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -6.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.85"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.28"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.199"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.55"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.833.01"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.932"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.147"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.79%  "

$ws.Range("E34").Value = "  -4.40%  "

$ws.Range("E35").Value = "  -8.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02712"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08717"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2436"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.937"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.46%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06792"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.82%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.92"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6878"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.296"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.34%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.62"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6351"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.29%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.251"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.904"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07720"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.49"
$ws.Range("D50").ClearFormats()

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.149"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.11%  "
